# #5: property boat&car done
#
# The "汽車" (car) sheet (3rd worksheet) had a broken header row 1 that was
# just a duplicate of the data row 2 instead of real column headers, and it
# was missing the trailing metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that every other
# property sheet in this workbook already has. This fixes the header row and
# appends the metadata columns to the data row, matching the pattern used on
# the other sheets (land/building/deposit/stock/fund/insurance/investment).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: replace the bogus duplicated-data header with real column names ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the new header cells (H1:N1) the same bold/bordered look as the
# existing header cells (B1:G1).
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2: keep the existing data (A2:G2 unchanged) and append the
# metadata columns present on every other sheet ---
# "date" must stay a plain text value ("2013-12-20"), not get auto-converted
# into a date serial number, so force a text number format before writing it.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-20"
$ws.Range("K2").Value = "尤美女"
$ws.Range("L2").Value = 1730
$ws.Range("M2").Value = "tmp84bd1"
$ws.Range("N2").Value = 46

# Match the plain formatting used by the rest of row 2 (B2:G2).
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
